# Update the "想去人数" (F column) counts across all sheets to the newly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3286
$ws.Range("F6").Value = 7685
$ws.Range("F9").Value = 686
$ws.Range("F12").Value = 1025
$ws.Range("F16").Value = 6105
$ws.Range("F17").Value = 55
$ws.Range("F21").Value = 1010
$ws.Range("F22").Value = 4233
$ws.Range("F24").Value = 313
$ws.Range("F25").Value = 128
$ws.Range("F26").Value = 1055
$ws.Range("F31").Value = 1041
$ws.Range("F37").Value = 196
$ws.Range("F42").Value = 1142
$ws.Range("F45").Value = 3149
$ws.Range("F46").Value = 85
$ws.Range("F47").Value = 382

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 259
$ws.Range("F23").Value = 28
$ws.Range("F25").Value = 6420
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 31

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 3038
$ws.Range("F6").Value = 1210
$ws.Range("F7").Value = 1278
$ws.Range("F9").Value = 532
$ws.Range("F11").Value = 8797
$ws.Range("F13").Value = 64

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3286
$ws.Range("F5").Value = 7685
$ws.Range("F6").Value = 1210
$ws.Range("F7").Value = 532
$ws.Range("F11").Value = 64
$ws.Range("F14").Value = 1025
$ws.Range("F17").Value = 259
$ws.Range("F19").Value = 6105
$ws.Range("F20").Value = 55
$ws.Range("F23").Value = 1010
$ws.Range("F25").Value = 128
$ws.Range("F26").Value = 1055
$ws.Range("F28").Value = 1041
$ws.Range("F33").Value = 196
$ws.Range("F39").Value = 28
$ws.Range("F40").Value = 3149
$ws.Range("F41").Value = 85
$ws.Range("F42").Value = 6420
